# Auto-generated edit script applying numeric corrections to Chocobo_Profits leve tables
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

$wsALC.Range("H17").Value = 948.2436
$wsALC.Range("J17").Value = 896.2381
$wsALC.Range("L17").Value = 2688.7143
$wsALC.Range("N17").Value = -3024.7143
$wsALC.Range("H33").Value = 187.35294
$wsALC.Range("J33").Value = 506.66666
$wsALC.Range("L33").Value = 506.66666
$wsALC.Range("N33").Value = -964.66666
$wsALC.Range("H76").Value = 3325
$wsALC.Range("I76").Value = 3100
$wsALC.Range("J76").Value = 3550
$wsALC.Range("K76").Value = 3100
$wsALC.Range("L76").Value = 3550
$wsALC.Range("M76").Value = -2785
$wsALC.Range("N76").Value = -4180
$wsALC.Range("H79").Value = 3325
$wsALC.Range("I79").Value = 3100
$wsALC.Range("J79").Value = 3550
$wsALC.Range("K79").Value = 3100
$wsALC.Range("L79").Value = 3550
$wsALC.Range("M79").Value = -2008
$wsALC.Range("N79").Value = -5734
$wsALC.Range("H135").Value = 725.8182
$wsALC.Range("I135").Value = 697.7143
$wsALC.Range("K135").Value = 6279.428699999999
$wsALC.Range("M135").Value = -3744.428699999999
$wsALC.Range("H138").Value = 5402.91
$wsALC.Range("I138").Value = 1276.7693
$wsALC.Range("J138").Value = 6019.46
$wsALC.Range("K138").Value = 3830.3079
$wsALC.Range("L138").Value = 18058.38
$wsALC.Range("M138").Value = 1309.6921
$wsALC.Range("N138").Value = -28338.38
$wsALC.Range("H141").Value = 30268.97
$wsALC.Range("J141").Value = 3178.4285
$wsALC.Range("L141").Value = 9535.2855
$wsALC.Range("N141").Value = -19895.2855
$wsARM.Range("H45").Value = 1110.3334
$wsARM.Range("I45").Value = 937.3333
$wsARM.Range("J45").Value = 1283.3334
$wsARM.Range("K45").Value = 937.3333
$wsARM.Range("L45").Value = 1283.3334
$wsARM.Range("M45").Value = -560.3333
$wsARM.Range("N45").Value = -2037.3334
$wsARM.Range("H122").Value = 2262.7778
$wsARM.Range("I122").Value = 1594.2941
$wsARM.Range("J122").Value = 3399.2
$wsARM.Range("K122").Value = 4782.8823
$wsARM.Range("L122").Value = 10197.6
$wsARM.Range("M122").Value = -2332.8823
$wsARM.Range("N122").Value = -15097.6
$wsARM.Range("H132").Value = 2450.432
$wsARM.Range("I132").Value = 1446.1936
$wsARM.Range("K132").Value = 4338.5808
$wsARM.Range("M132").Value = -1808.5808
$wsBSM.Range("H48").Value = 74800
$wsBSM.Range("J48").Value = 74800
$wsBSM.Range("L48").Value = 74800
$wsBSM.Range("N48").Value = -75630
$wsBSM.Range("H134").Value = 1886.6184
$wsBSM.Range("I134").Value = 1182.5555
$wsBSM.Range("K134").Value = 3547.6665
$wsBSM.Range("M134").Value = -1012.6665
$wsCRP.Range("H16").Value = 4832612
$wsCRP.Range("I16").Value = 8548365
$wsCRP.Range("J16").Value = 2133.5
$wsCRP.Range("K16").Value = 8548365
$wsCRP.Range("L16").Value = 2133.5
$wsCRP.Range("M16").Value = -8548078
$wsCRP.Range("N16").Value = -2707.5
$wsCRP.Range("H58").Value = 1382.17
$wsCRP.Range("I58").Value = 1530.0986
$wsCRP.Range("J58").Value = 1020
$wsCRP.Range("K58").Value = 1530.0986
$wsCRP.Range("L58").Value = 1020
$wsCRP.Range("M58").Value = -1327.0986
$wsCRP.Range("N58").Value = -1426
$wsCRP.Range("H98").Value = 45000
$wsCRP.Range("J98").Value = 45000
$wsCRP.Range("L98").Value = 45000
$wsCRP.Range("N98").Value = -49492
$wsCRP.Range("H106").Value = 29375
$wsCRP.Range("J106").Value = 29375
$wsCRP.Range("L106").Value = 29375
$wsCRP.Range("N106").Value = -31899
$wsCRP.Range("H109").Value = 34666.668
$wsCRP.Range("J109").Value = 34666.668
$wsCRP.Range("L109").Value = 34666.668
$wsCRP.Range("N109").Value = -36746.668
$wsCRP.Range("H113").Value = 4832612
$wsCRP.Range("I113").Value = 8548365
$wsCRP.Range("J113").Value = 2133.5
$wsCRP.Range("K113").Value = 8548365
$wsCRP.Range("L113").Value = 2133.5
$wsCRP.Range("M113").Value = -8546195
$wsCRP.Range("N113").Value = -6473.5
$wsCRP.Range("H134").Value = 4183.5
$wsCRP.Range("I134").Value = 5686.048
$wsCRP.Range("K134").Value = 17058.144
$wsCRP.Range("M134").Value = -14523.144
$wsCRP.Range("H136").Value = 1382.17
$wsCRP.Range("I136").Value = 1530.0986
$wsCRP.Range("J136").Value = 1020
$wsCRP.Range("K136").Value = 4590.2958
$wsCRP.Range("L136").Value = 3060
$wsCRP.Range("M136").Value = -2040.2958
$wsCRP.Range("N136").Value = -8160
$wsCRP.Range("H141").Value = 34933.332
$wsCRP.Range("J141").Value = 34933.332
$wsCRP.Range("L141").Value = 34933.332
$wsCRP.Range("N141").Value = -45293.332
$wsCUL.Range("H131").Value = 837.9041
$wsCUL.Range("I131").Value = 502.22223
$wsCUL.Range("J131").Value = 885.1094000000001
$wsCUL.Range("K131").Value = 1506.66669
$wsCUL.Range("L131").Value = 2655.3282
$wsCUL.Range("M131").Value = 3533.33331
$wsCUL.Range("N131").Value = -12735.3282
$wsGSM.Range("H43").Value = 16791.264
$wsGSM.Range("I43").Value = 1141.4286
$wsGSM.Range("J43").Value = 25920.334
$wsGSM.Range("K43").Value = 1141.4286
$wsGSM.Range("L43").Value = 25920.334
$wsGSM.Range("M43").Value = -990.4286
$wsGSM.Range("N43").Value = -26222.334
$wsGSM.Range("H46").Value = 30909.572
$wsGSM.Range("J46").Value = 30909.572
$wsGSM.Range("L46").Value = 30909.572
$wsGSM.Range("N46").Value = -31221.572
$wsGSM.Range("H57").Value = 38600
$wsGSM.Range("J57").Value = 38333.332
$wsGSM.Range("L57").Value = 38333.332
$wsGSM.Range("N57").Value = -39973.332
$wsGSM.Range("H70").Value = 6227.531
$wsGSM.Range("I70").Value = 5768.4
$wsGSM.Range("J70").Value = 7375.357
$wsGSM.Range("K70").Value = 5768.4
$wsGSM.Range("L70").Value = 7375.357
$wsGSM.Range("M70").Value = -5498.4
$wsGSM.Range("N70").Value = -7915.357
$wsGSM.Range("H73").Value = 6227.531
$wsGSM.Range("I73").Value = 5768.4
$wsGSM.Range("J73").Value = 7375.357
$wsGSM.Range("K73").Value = 5768.4
$wsGSM.Range("L73").Value = 7375.357
$wsGSM.Range("M73").Value = -4832.4
$wsGSM.Range("N73").Value = -9247.357
$wsGSM.Range("H80").Value = 22729756
$wsGSM.Range("I80").Value = 27780168
$wsGSM.Range("J80").Value = 2904.5
$wsGSM.Range("K80").Value = 27780168
$wsGSM.Range("L80").Value = 2904.5
$wsGSM.Range("M80").Value = -27779170
$wsGSM.Range("N80").Value = -4900.5
$wsGSM.Range("H83").Value = 22729756
$wsGSM.Range("I83").Value = 27780168
$wsGSM.Range("J83").Value = 2904.5
$wsGSM.Range("K83").Value = 138900840
$wsGSM.Range("L83").Value = 14522.5
$wsGSM.Range("M83").Value = -138895848
$wsGSM.Range("N83").Value = -24506.5
$wsGSM.Range("H112").Value = 0
$wsGSM.Range("J112").Value = 0
$wsGSM.Range("L112").Value = 0
$wsGSM.Range("N112").ClearContents()
$wsGSM.Range("H113").Value = 1376.0714
$wsGSM.Range("I113").Value = 978.25
$wsGSM.Range("J113").Value = 1906.5
$wsGSM.Range("K113").Value = 978.25
$wsGSM.Range("L113").Value = 1906.5
$wsGSM.Range("M113").Value = 1191.75
$wsGSM.Range("N113").Value = -6246.5
$wsGSM.Range("H122").Value = 3082.5715
$wsGSM.Range("I122").Value = 1434.5333
$wsGSM.Range("J122").Value = 7202.6665
$wsGSM.Range("K122").Value = 4303.5999
$wsGSM.Range("L122").Value = 21607.9995
$wsGSM.Range("M122").Value = -1853.5999
$wsGSM.Range("N122").Value = -26507.9995
$wsGSM.Range("H132").Value = 3679.8076
$wsGSM.Range("I132").Value = 2339.8667
$wsGSM.Range("K132").Value = 7019.6001
$wsGSM.Range("M132").Value = -4489.6001
$wsLTW.Range("H93").Value = 4832655.5
$wsLTW.Range("I93").Value = 10102333
$wsLTW.Range("J93").Value = 2117.25
$wsLTW.Range("K93").Value = 10102333
$wsLTW.Range("L93").Value = 2117.25
$wsLTW.Range("M93").Value = -10101085
$wsLTW.Range("N93").Value = -4613.25
$wsLTW.Range("H100").Value = 3000
$wsLTW.Range("I100").Value = 0
$wsLTW.Range("K100").Value = 0
$wsLTW.Range("M100").ClearContents()
$wsWVR.Range("H93").Value = 39750
$wsWVR.Range("J93").Value = 39750
$wsWVR.Range("L93").Value = 39750
$wsWVR.Range("N93").Value = -44742
$wsWVR.Range("H115").Value = 26574.074
$wsWVR.Range("J115").Value = 26574.074
$wsWVR.Range("L115").Value = 26574.074
$wsWVR.Range("N115").Value = -29708.074
